$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.128.53'
$ws.Range('E2').Value = '  -8.68%  '
$ws.Range('D3').Value = '3.177.51'
$ws.Range('E3').Value = '  -10.20%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '504.92'
$ws.Range('E5').Value = '  -9.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.27'
$ws.Range('E6').Value = '  -15.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.581'
$ws.Range('E7').Value = '  -10.16%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '3.175.12'
$ws.Range('E9').Value = '  -10.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.583'
$ws.Range('E10').Value = '  -12.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.70'
$ws.Range('E11').Value = '  -12.73%  '
$ws.Range('E12').Value = '  -11.99%  '
$ws.Range('E13').Value = '  -9.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.67'
$ws.Range('E14').Value = '  -12.98%  '
$ws.Range('D15').Value = '3.691.57'
$ws.Range('E15').Value = '  -9.77%  '
$ws.Range('D16').Value = '3.180.27'
$ws.Range('E16').Value = '  -9.87%  '
$ws.Range('D17').Value = '62.037.52'
$ws.Range('E17').Value = '  -8.34%  '
$ws.Range('E18').Value = '  -10.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.66'
$ws.Range('E19').Value = '  -9.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.54'
$ws.Range('E20').Value = '  -11.70%  '
$ws.Range('E21').Value = '  -10.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '354.51'
$ws.Range('E22').Value = '  -11.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.59'
$ws.Range('E23').Value = '  -10.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.67'
$ws.Range('E24').Value = '  -9.37%  '
$ws.Range('E25').Value = '  -1.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.48'
$ws.Range('E26').Value = '  -11.86%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.76'
$ws.Range('E27').Value = '  -2.79%  '
$ws.Range('E28').Value = '  -10.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.82'
$ws.Range('E29').Value = '  -12.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.93'
$ws.Range('E30').Value = '  -11.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '27.48'
$ws.Range('E31').Value = '  -12.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '607.76'
$ws.Range('E32').Value = '  -16.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.27'
$ws.Range('E33').Value = '  -11.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.79'
$ws.Range('E34').Value = '  -8.57%  '
$ws.Range('B35').Value = 'Dai'
$ws.Range('C35').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.100'
$ws.Range('E36').Value = '  -10.31%  '
$ws.Range('E37').Value = '  -13.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '35.37'
$ws.Range('E38').Value = '  -8.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.367'
$ws.Range('E39').Value = '  -7.28%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.997'
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0652'
$ws.Range('E41').Value = '  -5.76%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.794.35'
$ws.Range('E42').Value = '  -9.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.117'
$ws.Range('E43').Value = '  -10.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.33'
$ws.Range('E44').Value = '  -7.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.55'
$ws.Range('E45').Value = '  -7.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.54'
$ws.Range('E46').Value = '  -16.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0376'
$ws.Range('E47').Value = '  -8.28%  '
$ws.Range('E48').Value = '  -3.10%  '
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '131.83'
$ws.Range('E50').Value = '  -5.60%  '
$ws.Range('E51').Value = '  -10.53%  '
